$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A width (target stored width 15.42578125; the host quantizes
# ColumnWidth to a 1/6-character pixel grid, so 14.67 is the input that lands
# on the closest achievable stored width to the target)
$ws.Columns.Item(1).ColumnWidth = 14.67

# Update cell values A1:A33
$ws.Range("A1").Value = 0.072341782119728748
$ws.Range("A2").Value = -0.0059999999515270019
$ws.Range("A3").Value = 0.011248036374084336
$ws.Range("A4").Value = -0.0079999999312736492
$ws.Range("A5").Value = -0.0029999999707026603
$ws.Range("A6").Value = -0.0019999999788407052
$ws.Range("A7").Value = -0.0099999999148452368
$ws.Range("A8").Value = -0.009999999913109292
$ws.Range("A9").Value = -0.0019999999756792342
$ws.Range("A10").Value = -0.001999999975096145
$ws.Range("A11").Value = -0.0029999999671748157
$ws.Range("A12").Value = -0.0034999999632270296
$ws.Range("A13").Value = -0.0034999999641973645
$ws.Range("A14").Value = 0.0039371701147015159
$ws.Range("A15").Value = -0.00099999998503275123
$ws.Range("A16").Value = 0.010025298565358387
$ws.Range("A17").Value = -0.0019999999764612753
$ws.Range("A18").Value = -0.0039999999601896263
$ws.Range("A19").Value = -0.0039999999681108456
$ws.Range("A20").Value = -0.0039999999654227736
$ws.Range("A21").Value = -0.0039999999649529272
$ws.Range("A22").Value = -0.003999999964558576
$ws.Range("A23").Value = -0.0049999999529886097
$ws.Range("A24").Value = -0.019999999828372417
$ws.Range("A25").Value = -0.019999999825535575
$ws.Range("A26").Value = 0.072433875515187651
$ws.Range("A27").Value = -0.0024999999738439804
$ws.Range("A28").Value = -0.0019999999706801219
$ws.Range("A29").Value = -0.0069999999258518741
$ws.Range("A30").Value = -0.059999999502931001
$ws.Range("A31").Value = -0.0069999999283307801
$ws.Range("A32").Value = 0.039929749545276394
$ws.Range("A33").Value = -0.0039999999537432274
